$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 749.5
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 749.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 749.5
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -975.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8081.273
$ws.Range("I40").Value = 1199
$ws.Range("K40").Value = 1199
$ws.Range("M40").Value = -1024

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4687.4443
$ws.Range("J76").Value = 3500
$ws.Range("L76").Value = 3500
$ws.Range("N76").Value = -4130

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4687.4443
$ws.Range("J79").Value = 3500
$ws.Range("L79").Value = 3500
$ws.Range("N79").Value = -5684

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3293628
$ws.Range("I86").Value = 4602.5454
$ws.Range("J86").Value = 10529484
$ws.Range("K86").Value = 4602.5454
$ws.Range("L86").Value = 10529484
$ws.Range("M86").Value = -3479.5454
$ws.Range("N86").Value = -10531730

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3293628
$ws.Range("I89").Value = 4602.5454
$ws.Range("J89").Value = 10529484
$ws.Range("K89").Value = 23012.727
$ws.Range("L89").Value = 52647420
$ws.Range("M89").Value = -17396.727
$ws.Range("N89").Value = -52658652

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 135457.86
$ws.Range("I129").Value = 168530.58
$ws.Range("J129").Value = 3167
$ws.Range("K129").Value = 505591.74
$ws.Range("L129").Value = 9501
$ws.Range("M129").Value = -500591.74
$ws.Range("N129").Value = -19501

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1094.5454
$ws.Range("I132").Value = 971.3214
$ws.Range("K132").Value = 2913.9642
$ws.Range("M132").Value = -383.9642000000003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2786.6226
$ws.Range("I138").Value = 1711.6
$ws.Range("J138").Value = 2898.6042
$ws.Range("K138").Value = 5134.799999999999
$ws.Range("L138").Value = 8695.812600000001
$ws.Range("M138").Value = 5.200000000000728
$ws.Range("N138").Value = -18975.8126

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6794.9575
$ws.Range("I2").Value = 7997.1797
$ws.Range("K2").Value = 7997.1797
$ws.Range("M2").Value = -7884.1797

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3524.625
$ws.Range("I45").Value = 2195.6924
$ws.Range("K45").Value = 2195.6924
$ws.Range("M45").Value = -1818.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1667.76
$ws.Range("I74").Value = 1182.9
$ws.Range("K74").Value = 1182.9
$ws.Range("M74").Value = -308.9000000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1667.76
$ws.Range("I77").Value = 1182.9
$ws.Range("K77").Value = 5914.5
$ws.Range("M77").Value = -1546.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 6794.9575
$ws.Range("I116").Value = 7997.1797
$ws.Range("K116").Value = 7997.1797
$ws.Range("M116").Value = -5703.1797

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6794.9575
$ws.Range("I3").Value = 7997.1797
$ws.Range("K3").Value = 7997.1797
$ws.Range("M3").Value = -7883.1797

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 36574.17
$ws.Range("I20").Value = 2102.85
$ws.Range("J20").Value = 113177.11
$ws.Range("K20").Value = 2102.85
$ws.Range("L20").Value = 113177.11
$ws.Range("M20").Value = -1855.85
$ws.Range("N20").Value = -113671.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 4674.6665
$ws.Range("I37").Value = 24
$ws.Range("J37").Value = 7000
$ws.Range("K37").Value = 24
$ws.Range("L37").Value = 7000
$ws.Range("M37").Value = 113
$ws.Range("N37").Value = -7274

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 51333.332
$ws.Range("J92").Value = 51333.332
$ws.Range("L92").Value = 51333.332
$ws.Range("N92").Value = -56325.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 10529333
$ws.Range("I94").Value = 2381.4666
$ws.Range("K94").Value = 2381.4666
$ws.Range("M94").Value = -1930.4666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 12960.143
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4424690
$ws.Range("I122").Value = 10405746
$ws.Range("J122").Value = 3909.5217
$ws.Range("K122").Value = 31217238
$ws.Range("L122").Value = 11728.5651
$ws.Range("M122").Value = -31214788
$ws.Range("N122").Value = -16628.5651

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 143.5
$ws.Range("J2").Value = 236.5
$ws.Range("L2").Value = 1419
$ws.Range("N2").Value = -1645

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2923
$ws.Range("J36").Value = 3147.5
$ws.Range("L36").Value = 9442.5
$ws.Range("N36").Value = -9780.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 836.5357
$ws.Range("I38").Value = 83.29412000000001
$ws.Range("J38").Value = 2000.6364
$ws.Range("K38").Value = 249.88236
$ws.Range("L38").Value = 6001.9092
$ws.Range("M38").Value = 97.11763999999999
$ws.Range("N38").Value = -6695.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5759.294
$ws.Range("I122").Value = 3811.111
$ws.Range("J122").Value = 7951
$ws.Range("K122").Value = 11433.333
$ws.Range("L122").Value = 23853
$ws.Range("M122").Value = -8983.332999999999
$ws.Range("N122").Value = -28753

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4325.125
$ws.Range("I126").Value = 2932.3333
$ws.Range("J126").Value = 8503.5
$ws.Range("K126").Value = 8796.999899999999
$ws.Range("L126").Value = 25510.5
$ws.Range("M126").Value = -6326.999899999999
$ws.Range("N126").Value = -30450.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1431.24
$ws.Range("I55").Value = 1563.0667
$ws.Range("K55").Value = 1563.0667
$ws.Range("M55").Value = -1390.0667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4309.6665
$ws.Range("I122").Value = 3668.4546
$ws.Range("J122").Value = 7131
$ws.Range("K122").Value = 11005.3638
$ws.Range("L122").Value = 21393
$ws.Range("M122").Value = -8555.363799999999
$ws.Range("N122").Value = -26293

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10103444
$ws.Range("I81").Value = 2218.76
$ws.Range("J81").Value = 41669776
$ws.Range("K81").Value = 4437.52
$ws.Range("L81").Value = 83339552
$ws.Range("M81").Value = -3376.52
$ws.Range("N81").Value = -83341674

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 10103444
$ws.Range("I84").Value = 2218.76
$ws.Range("J84").Value = 41669776
$ws.Range("K84").Value = 22187.6
$ws.Range("L84").Value = 416697760
$ws.Range("M84").Value = -16883.6
$ws.Range("N84").Value = -416708368

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2411.7368
$ws.Range("I136").Value = 1887.4584
$ws.Range("K136").Value = 5662.3752
$ws.Range("M136").Value = -3112.3752
